# "Generate Report for Archive" - refresh localization status report.
#
# 1) Status text moves from "Ready for handoff" to "In Translation" on every
#    sheet that surfaces the per-file status (Overview!E:F, zh-cn!C, de-de!C).
# 2) The (now shorter) status text lets those Status columns be narrowed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- update status values -------------------------------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# --- narrow the Status columns to match the new, shorter text -------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
